$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

# Insert a new row above row 45 so the new translation entry
# ("common.selection.InventoryAromas.title") lands right after
# "common.selection.cancel.label" and before the other selection-title
# entries, shifting everything below it down by one row.
$ws.Rows.Item(45).Insert()

# Fill in the new row with the new translation key/value pair.
$ws.Cells.Item(45, 1).Value = "cs"
$ws.Cells.Item(45, 2).Value = "common.selection.InventoryAromas.title"
$ws.Cells.Item(45, 3).Value = "Výběr příchutě"

# Re-apply the remembered sort over the block that moved down with the
# insert, so the sheet's stored sort state follows the new row layout.
$ws.Sort.SortFields.Clear()
$key = $ws.Range("B46:B58")
$ws.Sort.SortFields.Add2($key, 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("A46:C58"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Restore the active selection to match the saved view state.
$ws.Range("B46").Select()
